$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 112
$ws.Cells.Item(2, 4).Value = 'UN'
$ws.Cells.Item(2, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(3, 3).Value = 78
$ws.Cells.Item(3, 4).Value = 'UN'
$ws.Cells.Item(3, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(4, 3).Value = 234
$ws.Cells.Item(4, 4).Value = 'UN'
$ws.Cells.Item(4, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(5, 1).Value = '10.937.338'
$ws.Cells.Item(5, 2).Value = 'Cadeado 25x26x17x4,5x14x13mm capa termop'
$ws.Cells.Item(5, 4).Value = 'UN'
$ws.Cells.Item(5, 5).Value = 'ZD'
$ws.Cells.Item(5, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(6, 1).Value = '10.544.315'
$ws.Cells.Item(6, 2).Value = 'Cadeado 25x26x17x4,5x14x13mm capa termop'
$ws.Cells.Item(6, 4).Value = 'UN'
$ws.Cells.Item(6, 5).Value = 'ZD'
$ws.Cells.Item(6, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(7, 1).Value = '10.667.890'
$ws.Cells.Item(7, 2).Value = 'Cadeado 30x30x20x5x17x13,5mm capa termop'
$ws.Cells.Item(7, 3).Value = 403
$ws.Cells.Item(7, 5).Value = 'VB'
$ws.Cells.Item(7, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(8, 1).Value = '10.394.315'
$ws.Cells.Item(8, 2).Value = 'Cadeado 30x30x20x5x17x13,5mm capa termop'
$ws.Cells.Item(8, 3).Value = 186
$ws.Cells.Item(8, 5).Value = 'ZS'
$ws.Cells.Item(8, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(9, 3).Value = 7
$ws.Cells.Item(9, 4).Value = 'UN'
$ws.Cells.Item(9, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(10, 3).Value = 47
$ws.Cells.Item(10, 4).Value = 'UN'
$ws.Cells.Item(10, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(11, 2).Value = 'Chave p/volante válvula'
$ws.Cells.Item(11, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 'UN'
$ws.Cells.Item(12, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(13, 3).Value = 6
$ws.Cells.Item(13, 4).Value = 'UN'
$ws.Cells.Item(13, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(14, 3).Value = 4
$ws.Cells.Item(14, 4).Value = 'UN'
$ws.Cells.Item(14, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(15, 3).Value = 11
$ws.Cells.Item(15, 4).Value = 'UN'
$ws.Cells.Item(15, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(16, 2).Value = 'Cola a frio borr. sintética bisn 85ml'
$ws.Cells.Item(16, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(17, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(18, 3).Value = 34
$ws.Cells.Item(18, 4).Value = 'UN'
$ws.Cells.Item(18, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(19, 3).Value = 21
$ws.Cells.Item(19, 4).Value = 'UN'
$ws.Cells.Item(19, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(20, 3).Value = 599
$ws.Cells.Item(20, 4).Value = 'UN'
$ws.Cells.Item(20, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(21, 3).Value = 1000
$ws.Cells.Item(21, 4).Value = 'UN'
$ws.Cells.Item(21, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(22, 3).Value = 2300
$ws.Cells.Item(22, 4).Value = 'UN'
$ws.Cells.Item(22, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(23, 3).Value = 800
$ws.Cells.Item(23, 4).Value = 'UN'
$ws.Cells.Item(23, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(24, 3).Value = 1000
$ws.Cells.Item(24, 4).Value = 'UN'
$ws.Cells.Item(24, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(25, 1).Value = '11.331.124'
$ws.Cells.Item(25, 2).Value = 'Extintor incên. cil.AC 6-A:40-B:C'
$ws.Cells.Item(25, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(26, 1).Value = '10.026.438'
$ws.Cells.Item(26, 2).Value = 'Extintor incên. cil.AI 40-B:C'
$ws.Cells.Item(26, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(27, 1).Value = '11.809.118'
$ws.Cells.Item(27, 2).Value = 'Isobutileno 100ppm cil.0,1m³'
$ws.Cells.Item(27, 4).Value = 'M3'
$ws.Cells.Item(27, 6).Value = 'Equipamentos de Medição e Gases'
$ws.Cells.Item(27, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(28, 1).Value = '11.874.356'
$ws.Cells.Item(28, 2).Value = 'Kit absorvedor p/derramamento'
$ws.Cells.Item(28, 3).Value = 3
$ws.Cells.Item(28, 4).Value = 'UN'
$ws.Cells.Item(28, 6).Value = 'Equipamentos de Emergência Ambiental'
$ws.Cells.Item(28, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(29, 1).Value = '12.088.160'
$ws.Cells.Item(29, 2).Value = 'LGE(líquido gerador de espuma) 1-3% 1000'
$ws.Cells.Item(29, 3).Value = 10000
$ws.Cells.Item(29, 4).Value = 'L'
$ws.Cells.Item(29, 5).Value = 'ZS'
$ws.Cells.Item(29, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(30, 1).Value = '10.734.499'
$ws.Cells.Item(30, 2).Value = 'LGE(líquido gerador de espuma) 3% 1000L'
$ws.Cells.Item(30, 3).Value = 13000
$ws.Cells.Item(30, 4).Value = 'L'
$ws.Cells.Item(30, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(31, 1).Value = '10.734.501'
$ws.Cells.Item(31, 2).Value = 'LGE(líquido gerador de espuma) 3-6% 1000'
$ws.Cells.Item(31, 3).Value = 14000
$ws.Cells.Item(31, 4).Value = 'L'
$ws.Cells.Item(31, 5).Value = 'ZD'
$ws.Cells.Item(31, 6).Value = 'Equipamentos de Combate à Incêndio'
$ws.Cells.Item(31, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(32, 1).Value = '10.997.331'
$ws.Cells.Item(32, 2).Value = 'Lacre plást. 30cm polipropil.'
$ws.Cells.Item(32, 3).Value = 2500
$ws.Cells.Item(32, 4).Value = 'UN'
$ws.Cells.Item(32, 5).Value = 'ZS'
$ws.Cells.Item(32, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(33, 1).Value = '10.994.961'
$ws.Cells.Item(33, 2).Value = 'Lacre plást. D.1,5x500mm plástico'
$ws.Cells.Item(33, 3).Value = 4400
$ws.Cells.Item(33, 5).Value = 'VB'
$ws.Cells.Item(33, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(34, 1).Value = '11.158.686'
$ws.Cells.Item(34, 2).Value = 'Lacre plást. cordoalha AC D.1,5x600mm pl'
$ws.Cells.Item(34, 6).Value = 'Equipamentos de Libra'
$ws.Cells.Item(34, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(35, 1).Value = '13.267.594'
$ws.Cells.Item(35, 2).Value = 'Lanterna cap. 3x LR03(AAA) 1,5V Ex ia II'
$ws.Cells.Item(35, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(36, 1).Value = '13.445.279'
$ws.Cells.Item(36, 2).Value = 'Livreto 85x130mm 93fls. (AST)'
$ws.Cells.Item(36, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(37, 1).Value = '12.775.577'
$ws.Cells.Item(37, 2).Value = 'Livreto padrão básico segurança 85x130mm (PBS)'
$ws.Cells.Item(37, 4).Value = 'UN'
$ws.Cells.Item(37, 5).Value = 'ZD'
$ws.Cells.Item(37, 6).Value = 'Outros'
$ws.Cells.Item(37, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(38, 1).Value = '10.178.332'
$ws.Cells.Item(38, 2).Value = 'Lona plást. PRT 6,0mx 170micra x x50m'
$ws.Cells.Item(38, 3).Value = 50
$ws.Cells.Item(38, 4).Value = 'M'
$ws.Cells.Item(38, 5).Value = 'ZB'
$ws.Cells.Item(38, 6).Value = 'Equipamentos de Emergência Ambiental'
$ws.Cells.Item(38, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(39, 1).Value = '10.494.563'
$ws.Cells.Item(39, 2).Value = 'Lubrificante em aerosol fr 300ml'
$ws.Cells.Item(39, 3).Value = 146
$ws.Cells.Item(39, 4).Value = 'UN'
$ws.Cells.Item(39, 5).Value = 'VB'
$ws.Cells.Item(39, 6).Value = 'Outros'
$ws.Cells.Item(39, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(40, 1).Value = '10.057.533'
$ws.Cells.Item(40, 2).Value = 'Mangueira incêndio DN 40mm NBR11861 tipo'
$ws.Cells.Item(40, 3).Value = 50
$ws.Cells.Item(40, 4).Value = 'UN'
$ws.Cells.Item(40, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(41, 1).Value = '10.426.638'
$ws.Cells.Item(41, 2).Value = 'Mangueira incêndio DN 65mm NBR11861 tipo'
$ws.Cells.Item(41, 3).Value = 70
$ws.Cells.Item(41, 4).Value = 'UN'
$ws.Cells.Item(41, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(42, 1).Value = '10.203.856'
$ws.Cells.Item(42, 2).Value = 'Marcador esfer. p/metal pta.3mm amarela'
$ws.Cells.Item(42, 3).Value = 9
$ws.Cells.Item(42, 4).Value = 'UN'
$ws.Cells.Item(42, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(43, 1).Value = '11.006.532'
$ws.Cells.Item(43, 2).Value = 'Marcador esfer. p/metal pta.3mm vermelha'
$ws.Cells.Item(43, 3).Value = 16
$ws.Cells.Item(43, 5).Value = 'ZD'
$ws.Cells.Item(43, 6).Value = 'Outros'
$ws.Cells.Item(43, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(44, 2).Value = 'Mistura gasosa cil.34L 10ppm'
$ws.Cells.Item(44, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(45, 2).Value = 'Mistura gasosa cil.58L 1,45% 60ppm 15% 2'
$ws.Cells.Item(45, 4).Value = 'M3'
$ws.Cells.Item(45, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(46, 2).Value = 'Protetor de FDS'
$ws.Cells.Item(46, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(47, 3).Value = 16
$ws.Cells.Item(47, 4).Value = 'UN'
$ws.Cells.Item(47, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(48, 3).Value = 7
$ws.Cells.Item(48, 4).Value = 'UN'
$ws.Cells.Item(48, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(49, 3).Value = 2
$ws.Cells.Item(49, 4).Value = 'UN'
$ws.Cells.Item(49, 7).Value = '09/08/2025 09:01:03'
$ws.Cells.Item(50, 3).Value = 2
$ws.Cells.Item(50, 4).Value = 'UN'
$ws.Cells.Item(50, 7).Value = '09/08/2025 09:01:03'
